# Applies the "D6_Protestant Ethic" relabeling + expansion of the
# en-dash-range "Items Used" shorthand into explicit, comma-separated
# item lists on the Schema sheet, then restores the Items sheet's
# scroll position and leaves Schema active with C28 selected.

$wb = $excel.ActiveWorkbook

$schema = $wb.Worksheets.Item("Schema")

# Row 23 ("Protestant Ethic" indicator) was renumbered from D5_ to D6_
$schema.Range("A23").Value = "D6_Protestant Ethic"

# "Items Used" column: replace the en-dash shorthand ranges with the
# fully spelled-out, comma-separated item codes.
$schema.Range("C2").Value  = "F115, F116, F117"
$schema.Range("C3").Value  = "E025, e026, E027, E028,E029"
$schema.Range("C5").Value  = "A065, A066, A067, A068, A069, A070, A070, A071, A072, A073, A074"
$schema.Range("C7").Value  = "E069_01, E069_02, E069_04, E069_05, E069_06, E069_07, E069_08,  E069_17"
$schema.Range("C9").Value  = "E114, E115, E116"
$schema.Range("C11").Value = "E154, E155, E157, E158"
$schema.Range("C12").Value = "E159, E160, E161, E162"
$schema.Range("C21").Value = "E120, E121, E122, E123"
$schema.Range("C22").Value = "E035, E036, E037, E038, E039"
$schema.Range("C23").Value = "C036, C037, C038, C039"

# Reset the Items sheet's scroll/selection state (it was previously
# scrolled down to row 61), then re-activate Schema with C28 selected
# so it remains the tab shown on open.
$items = $wb.Worksheets.Item("Items")
$items.Activate()

$schema.Activate()
$schema.Range("C28").Select() | Out-Null
